# Automatische test-sync: 2025-07-22 17:50:50
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 21
$ws.Cells.Item($row, 1).Value = "Wat zijn jullie openingstijden?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #1: Wat zijn jullie openingstijden?"
$ws.Cells.Item($row, 4).Value = "Openingstijden / Locatie"
$ws.Cells.Item($row, 5).Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$ws.Cells.Item($row, 6).Value = "2025-07-22 17:50:01"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Ja"

# Extend conditional-formatting ranges to cover the newly added row
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $ws.Range($col + "2:" + $col + "20")
    $newRange = $ws.Range($col + "2:" + $col + "21")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Openingstijden / Locatie"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 3
